$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I29").Value = 9
$ws.Range("H29").Value = 2193
$ws.Range("M29").Value = 254
$ws.Range("J29").Value = 5833
$ws.Range("L29").Value = 17499
$ws.Range("N29").Value = -18061
$ws.Range("K29").Value = 27
$ws.Range("L113").Value = 4933
$ws.Range("H113").Value = 4933
$ws.Range("J113").Value = 4933
$ws.Range("N113").Value = -11441
$ws.Range("H125").Value = 185009.88
$ws.Range("J125").Value = 209291.2
$ws.Range("L125").Value = 1883620.8
$ws.Range("N125").Value = -1888540.8
$ws.Range("H127").Value = 1900
$ws.Range("K127").Value = 5700
$ws.Range("I127").Value = 1900
$ws.Range("M127").Value = -740
$ws.Range("M132").Value = -5431.2338
$ws.Range("K132").Value = 7961.2338
$ws.Range("H132").Value = 2765.125
$ws.Range("I132").Value = 2653.7446
$ws.Range("M138").Value = -5346.5
$ws.Range("I138").Value = 3495.5
$ws.Range("K138").Value = 10486.5
$ws.Range("H138").Value = 4852.6904
$ws.Range("I141").Value = 1850.4
$ws.Range("K141").Value = 5551.200000000001
$ws.Range("M141").Value = -371.2000000000007
$ws.Range("H141").Value = 4731.8887

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K2").Value = 1333.6
$ws.Range("M2").Value = -1220.6
$ws.Range("I2").Value = 1333.6
$ws.Range("H2").Value = 1246
$ws.Range("K61").Value = 2250.7334
$ws.Range("H61").Value = 20838114
$ws.Range("M61").Value = -2038.7334
$ws.Range("I61").Value = 2250.7334
$ws.Range("H74").Value = 3859.6
$ws.Range("N74").Value = -8153.6665
$ws.Range("M74").Value = -1288.2222
$ws.Range("L74").Value = 6405.6665
$ws.Range("K74").Value = 2162.2222
$ws.Range("I74").Value = 2162.2222
$ws.Range("J74").Value = 6405.6665
$ws.Range("H77").Value = 3859.6
$ws.Range("L77").Value = 32028.3325
$ws.Range("I77").Value = 2162.2222
$ws.Range("N77").Value = -40764.3325
$ws.Range("K77").Value = 10811.111
$ws.Range("J77").Value = 6405.6665
$ws.Range("M77").Value = -6443.111000000001
$ws.Range("H88").Value = 9500
$ws.Range("L88").Value = 9500
$ws.Range("N88").Value = -10312
$ws.Range("J88").Value = 9500
$ws.Range("J91").Value = 9500
$ws.Range("N91").Value = -12308
$ws.Range("L91").Value = 9500
$ws.Range("H91").Value = 9500
$ws.Range("M116").Value = 960.4000000000001
$ws.Range("K116").Value = 1333.6
$ws.Range("H116").Value = 1246
$ws.Range("I116").Value = 1333.6
$ws.Range("K136").Value = 6752.2002
$ws.Range("M136").Value = -4202.2002
$ws.Range("H136").Value = 20838114
$ws.Range("I136").Value = 2250.7334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 1333.6
$ws.Range("M3").Value = -1219.6
$ws.Range("K3").Value = 1333.6
$ws.Range("H3").Value = 1246
$ws.Range("M99").Value = 221.7778000000001
$ws.Range("I99").Value = 1276.2222
$ws.Range("K99").Value = 1276.2222
$ws.Range("H99").Value = 1845.9231
$ws.Range("M134").Value = -4215
$ws.Range("N134").Value = -36849.60000000001
$ws.Range("I134").Value = 2250
$ws.Range("L134").Value = 31779.6
$ws.Range("J134").Value = 10593.2
$ws.Range("H134").Value = 3553.625
$ws.Range("K134").Value = 6750

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1500
$ws.Range("L31").Value = 10363.883
$ws.Range("J31").Value = 10363.883
$ws.Range("H31").Value = 4986.433
$ws.Range("I31").Value = 2860.465
$ws.Range("K31").Value = 2860.465
$ws.Range("M31").Value = -2565.465
$ws.Range("N31").Value = -10953.883
$ws.Range("N34").Value = -10767.883
$ws.Range("H34").Value = 4986.433
$ws.Range("K34").Value = 2860.465
$ws.Range("L34").Value = 10363.883
$ws.Range("M34").Value = -2658.465
$ws.Range("I34").Value = 2860.465
$ws.Range("J34").Value = 10363.883
$ws.Range("K58").Value = 3448.25
$ws.Range("M58").Value = -3245.25
$ws.Range("H58").Value = 7144.5264
$ws.Range("I58").Value = 3448.25
$ws.Range("H74").Value = 179874.5
$ws.Range("M74").ClearContents()
$ws.Range("K74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("H77").Value = 179874.5
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("M99").Value = -52.23070000000007
$ws.Range("I99").Value = 1550.2307
$ws.Range("K99").Value = 1550.2307
$ws.Range("H99").Value = 1744.098
$ws.Range("H113").Value = 1500
$ws.Range("M117").ClearContents()
$ws.Range("H117").Value = 199950
$ws.Range("K117").Value = 0
$ws.Range("J117").Value = 199950
$ws.Range("I117").Value = 0
$ws.Range("L117").Value = 199950
$ws.Range("N117").Value = -209128
$ws.Range("K122").Value = 3847.9998
$ws.Range("M122").Value = -1397.9998
$ws.Range("H122").Value = 1278.8235
$ws.Range("I122").Value = 1282.6666
$ws.Range("I126").Value = 1550.2307
$ws.Range("M126").Value = -2180.6921
$ws.Range("K126").Value = 4650.6921
$ws.Range("H126").Value = 1744.098
$ws.Range("M134").Value = -9351.375
$ws.Range("I134").Value = 3962.125
$ws.Range("H134").Value = 5521.8887
$ws.Range("K134").Value = 11886.375
$ws.Range("K136").Value = 10344.75
$ws.Range("M136").Value = -7794.75
$ws.Range("H136").Value = 7144.5264
$ws.Range("I136").Value = 3448.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I40").Value = 1193.6111
$ws.Range("K40").Value = 4774.4444
$ws.Range("M40").Value = -4705.4444
$ws.Range("H40").Value = 1183.421
$ws.Range("K64").Value = 12000
$ws.Range("H64").Value = 6799.8
$ws.Range("M64").Value = -11730
$ws.Range("I64").Value = 4000
$ws.Range("H67").Value = 6799.8
$ws.Range("K67").Value = 12000
$ws.Range("I67").Value = 4000
$ws.Range("M67").Value = -11064
$ws.Range("N69").ClearContents()
$ws.Range("H69").Value = 482
$ws.Range("L69").Value = 0
$ws.Range("K69").Value = 1446
$ws.Range("I69").Value = 482
$ws.Range("J69").Value = 0
$ws.Range("M69").Value = -635
$ws.Range("J72").Value = 0
$ws.Range("H72").Value = 482
$ws.Range("M72").Value = -282
$ws.Range("K72").Value = 4338
$ws.Range("N72").ClearContents()
$ws.Range("I72").Value = 482
$ws.Range("L72").Value = 0
$ws.Range("N120").Value = -99667
$ws.Range("M120").Value = -24240.571
$ws.Range("K120").Value = 29078.571
$ws.Range("L120").Value = 89991
$ws.Range("I120").Value = 9692.857
$ws.Range("J120").Value = 29997
$ws.Range("H120").Value = 14204.889
$ws.Range("I126").Value = 3999
$ws.Range("N126").Value = -83806
$ws.Range("L126").Value = 73926
$ws.Range("J126").Value = 24642
$ws.Range("M126").Value = -7057
$ws.Range("K126").Value = 11997
$ws.Range("H126").Value = 22061.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 605
$ws.Range("K97").Value = 395.5
$ws.Range("I97").Value = 395.5
$ws.Range("M97").Value = 100.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J7").Value = 4999
$ws.Range("L7").Value = 4999
$ws.Range("N7").Value = -5223
$ws.Range("H7").Value = 10782.5
$ws.Range("M55").Value = -226.8421
$ws.Range("N55").Value = -1344.4211
$ws.Range("L55").Value = 998.4211
$ws.Range("I55").Value = 399.8421
$ws.Range("J55").Value = 998.4211
$ws.Range("H55").Value = 699.1316
$ws.Range("K55").Value = 399.8421
$ws.Range("L100").Value = 2199.5
$ws.Range("N100").Value = -3281.5
$ws.Range("J100").Value = 2199.5
$ws.Range("H100").Value = 4689.2144
$ws.Range("N126").Value = -19937
$ws.Range("L126").Value = 14997
$ws.Range("J126").Value = 4999
$ws.Range("H126").Value = 10782.5
$ws.Range("N130").Value = -20469
$ws.Range("H130").Value = 10429
$ws.Range("L130").Value = 10429
$ws.Range("J130").Value = 10429
$ws.Range("K136").Value = 14179.8939
$ws.Range("M136").Value = -11629.8939
$ws.Range("I136").Value = 4726.6313

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K20").Value = 49800
$ws.Range("N20").Value = -17480
$ws.Range("M20").Value = -49560
$ws.Range("L20").Value = 17000
$ws.Range("J20").Value = 17000
$ws.Range("I20").Value = 49800
$ws.Range("H20").Value = 33400
$ws.Range("M132").Value = -6511.470499999999
$ws.Range("K132").Value = 9041.470499999999
$ws.Range("H132").Value = 3882.125
$ws.Range("I132").Value = 3013.8235
$ws.Range("K136").Value = 11174.2038
$ws.Range("M136").Value = -8624.203799999999
$ws.Range("H136").Value = 4574.681
$ws.Range("L136").Value = 19971.15
$ws.Range("I136").Value = 3724.7346
$ws.Range("N136").Value = -25071.15
$ws.Range("J136").Value = 6657.05

Write-Output "Applied all changes"